$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:M1): bold font, thin box border, centered horizontal / top vertical alignment
$headerRange = $ws.Range("A1:M1")
$headerRange.Font.Bold = $true
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Refresh row-2 usage counters (Groq / llama-3.1-8b-instant)
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 0.0008333333333333334
$ws.Range("K2").Value = 3213
$ws.Range("L2").Value = 0.006426
